$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data: pepino is "na_terra" (column B = "x")
$ws.Range("A14").Value = "pepino"
$ws.Range("B14").Value = "x"

# Move selection to A15 (next empty row), matching post-edit saved selection
$ws.Range("A15").Select()
